$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.377.33'
$ws.Range("E2").Value = '  +0.02%  '
$ws.Range("D3").Value = '1.874.74'
$ws.Range("E3").Value = '  -0.84%  '
$ws.Range("D4").Value = '''0.9999'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''238.41'
$ws.Range("E5").Value = '  +0.37%  '
$ws.Range("D6").Value = '''0.9999'
$ws.Range("D7").Value = '''0.4789'
$ws.Range("E7").Value = '  -1.12%  '
$ws.Range("D8").Value = '''0.2832'
$ws.Range("E8").Value = '  -2.34%  '
$ws.Range("D9").Value = '''0.06522'
$ws.Range("E9").Value = '  -1.22%  '
$ws.Range("D10").Value = '1.878.83'
$ws.Range("E10").Value = '  -0.56%  '
$ws.Range("D11").Value = '''0.07457'
$ws.Range("E11").Value = '  +1.68%  '
$ws.Range("D12").Value = '''16.61'
$ws.Range("E12").Value = '  -1.81%  '
$ws.Range("D13").Value = '''5.098'
$ws.Range("E13").Value = '  -1.22%  '
$ws.Range("D14").Value = '''88.28'
$ws.Range("E14").Value = '  +0.74%  '
$ws.Range("D15").Value = '''0.6582'
$ws.Range("E15").Value = '  -0.52%  '
$ws.Range("D16").Value = '30.353.31'
$ws.Range("E16").Value = '  +0.14%  '
$ws.Range("D17").Value = '''13.33'
$ws.Range("E17").Value = '  -0.76%  '
$ws.Range("D18").Value = '''0.9997'
$ws.Range("E18").Value = '  -0.04%  '
$ws.Range("D19").Value = '''0.000007601'
$ws.Range("E19").Value = '  -2.25%  '
$ws.Range("D20").Value = '2.115.41'
$ws.Range("E20").Value = '  -1.10%  '
$ws.Range("D21").Value = '''5.312'
$ws.Range("E21").Value = '  -2.17%  '
$ws.Range("E22").Value = '  +0.07%  '
$ws.Range("D23").Value = '''219.20'
$ws.Range("E23").Value = '  +12.86%  '
$ws.Range("D24").Value = '''6.228'
$ws.Range("E24").Value = '  +0.71%  '
$ws.Range("D25").Value = '''9.348'
$ws.Range("E25").Value = '  -0.13%  '
$ws.Range("D26").Value = '''167.76'
$ws.Range("E26").Value = '  +2.12%  '
$ws.Range("D27").Value = '''18.47'
$ws.Range("E27").Value = '  +1.55%  '
$ws.Range("D28").Value = '''1.978'
$ws.Range("E28").Value = '  +2.21%  '
$ws.Range("D29").Value = '''1.460'
$ws.Range("E29").Value = '  +0.40%  '
$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").Value = '''0.09428'
$ws.Range("E30").Value = '  +3.12%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = '''4.316'
$ws.Range("E31").Value = '  +0.28%  '
$ws.Range("D32").Value = '''4.043'
$ws.Range("E32").Value = '  +0.03%  '
$ws.Range("D33").Value = '''0.05065'
$ws.Range("E33").Value = '  -0.39%  '
$ws.Range("D34").Value = '''1.206'
$ws.Range("E34").Value = '  +5.83%  '
$ws.Range("D35").Value = '''0.7512'
$ws.Range("E35").Value = '  +2.85%  '
$ws.Range("D36").Value = '''2.711'
$ws.Range("E36").Value = '  +0.20%  '
$ws.Range("D37").Value = '''0.01825'
$ws.Range("E37").Value = '  +1.95%  '
$ws.Range("E38").Value = '  -1.28%  '
$ws.Range("D39").Value = '''2.068'
$ws.Range("E39").Value = '  -0.16%  '
$ws.Range("D40").Value = '''0.9051'
$ws.Range("E40").Value = '  -1.67%  '
$ws.Range("D41").Value = '''106.91'
$ws.Range("E41").Value = '  +1.08%  '
$ws.Range("D42").Value = '''5.893'
$ws.Range("E42").Value = '  +0.12%  '
$ws.Range("D43").Value = '''0.4280'
$ws.Range("E43").Value = '  -0.68%  '
$ws.Range("E44").Value = '  +0.14%  '
$ws.Range("D45").Value = '''7.389'
$ws.Range("E45").Value = '  -1.32%  '
$ws.Range("D46").Value = '''64.53'
$ws.Range("E46").Value = '  -0.55%  '
$ws.Range("E47").Value = '  -3.60%  '
$ws.Range("D48").Value = '''1.477'
$ws.Range("E48").Value = '  -7.55%  '
$ws.Range("D49").Value = '''8.966'
$ws.Range("E49").Value = '  -0.10%  '
$ws.Range("D50").Value = '''33.61'
$ws.Range("E50").Value = '  -0.95%  '
$ws.Range("D51").Value = '''0.3898'
$ws.Range("E51").Value = '  +1.18%  '
